# Regenerate save_data: replace column G ("K") values for rows 2-48
# with newly computed/recalculated values (per commit: "regen save_data
# to use K instead of Strike#, regen std/mean, calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..48 (column G), in row order.
$kValues = @(
    1, 2, 1, 1, 5, 1, 3, 5, 1, 2,
    1, 4, 1, 3, 2, 3, 1, 0, 0, 0,
    0, 2, 0, 2, 0, 3, 2, 1, 1, 2,
    0, 2, 2, 0, 0, 6, 3, 4, 2, 1,
    1, 0, 0, 1, 1, 0, 2
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
